$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename E1/F1/G1 ---
$ws.Range("E1").Value = "self_piece"
$ws.Range("F1").Value = "win_tour"
$ws.Range("G1").Value = "peer_piece"

# --- Data updates (columns E/F/G, rows 2-7) ---
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1

$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0

$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1

$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0

$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1

# --- Column G width ---
$ws.Columns.Item(7).ColumnWidth = 11.14

# --- Selection / active cell ---
$ws.Range("F11").Select() | Out-Null

# --- Window size (app window, mirrors the workbookView bookViews entry) ---
$excel.ActiveWindow.Width = 12130
$excel.ActiveWindow.Height = 12060
